# Auto-generated script applying the 2023-03-30 daily violent-crime data update.
# For every affected worksheet, the year-to-date (2023 / column J) totals are
# incremented by the counts of the newly-added day, and a couple of historical
# (column B / 2015) corrections are applied as recorded in the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 1559
$ws.Range('J3').Value = 1633
$ws.Range('B4').Value = 1670
$ws.Range('J4').Value = 368
$ws.Range('J5').Value = 113
$ws.Range('J6').Value = 2132
$ws.Range('B7').Value = 23302
$ws.Range('J7').Value = 5805

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J2').Value = 23
$ws.Range('J7').Value = 66

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J3').Value = 67
$ws.Range('J7').Value = 195

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J2').Value = 42
$ws.Range('J3').Value = 85
$ws.Range('J7').Value = 205

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J3').Value = 8
$ws.Range('J6').Value = 8
$ws.Range('J7').Value = 39

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('J6').Value = 24
$ws.Range('J7').Value = 59

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J2').Value = 40
$ws.Range('J7').Value = 147

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J5').Value = 16
$ws.Range('J6').Value = 58
$ws.Range('J7').Value = 165
$ws.Range('J8').Value = 359
$ws.Range('J10').Value = 35
$ws.Range('J11').Value = 71
$ws.Range('J19').Value = 203
$ws.Range('J20').Value = 121
$ws.Range('J26').Value = 8
$ws.Range('J29').Value = 322
$ws.Range('J31').Value = 39
$ws.Range('J33').Value = 242
$ws.Range('J37').Value = 195
$ws.Range('J41').Value = 37
$ws.Range('J42').Value = 227
$ws.Range('J48').Value = 46
$ws.Range('J49').Value = 36
$ws.Range('J52').Value = 130
$ws.Range('J53').Value = 54
$ws.Range('J54').Value = 110
$ws.Range('J60').Value = 33
$ws.Range('B63').Value = 374
$ws.Range('J63').Value = 25
$ws.Range('J65').Value = 147
$ws.Range('J67').Value = 205
$ws.Range('J71').Value = 28
$ws.Range('J73').Value = 54
$ws.Range('J76').Value = 89
$ws.Range('J77').Value = 43
$ws.Range('J78').Value = 76
$ws.Range('J79').Value = 181
$ws.Range('J80').Value = 14
$ws.Range('J83').Value = 142
$ws.Range('J84').Value = 59
$ws.Range('J85').Value = 263
$ws.Range('J86').Value = 31
$ws.Range('J88').Value = 60
$ws.Range('J89').Value = 66
$ws.Range('J94').Value = 49
$ws.Range('J95').Value = 85
$ws.Range('J98').Value = 39
$ws.Range('J100').Value = 12
$ws.Range('B101').Value = 23302
$ws.Range('J101').Value = 5805

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J6').Value = 44
$ws.Range('J7').Value = 142

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('J2').Value = 30
$ws.Range('J7').Value = 85

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 61
$ws.Range('J6').Value = 94
$ws.Range('J7').Value = 242

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('J4').Value = 3
$ws.Range('J7').Value = 36

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J2').Value = 32
$ws.Range('J7').Value = 110

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J3').Value = 119
$ws.Range('J4').Value = 15
$ws.Range('J7').Value = 322

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J6').Value = 74
$ws.Range('J7').Value = 203

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J3').Value = 6
$ws.Range('J7').Value = 46

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J6').Value = 52
$ws.Range('J7').Value = 89

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J2').Value = 66
$ws.Range('J7').Value = 263

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J6').Value = 22
$ws.Range('J7').Value = 58

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('J3').Value = 7
$ws.Range('J7').Value = 37

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J2').Value = 48
$ws.Range('J6').Value = 123
$ws.Range('J7').Value = 227

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('J6').Value = 14
$ws.Range('J7').Value = 35

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('J2').Value = 19
$ws.Range('J7').Value = 76

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J3').Value = 65
$ws.Range('J6').Value = 53
$ws.Range('J7').Value = 181

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J2').Value = 32
$ws.Range('J5').Value = 3
$ws.Range('J6').Value = 35
$ws.Range('J7').Value = 121

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range('J5').Value = 5
$ws.Range('J6').Value = 12

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J3').Value = 42
$ws.Range('J7').Value = 130

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('J6').Value = 32
$ws.Range('J7').Value = 49

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('J6').Value = 22
$ws.Range('J7').Value = 39

$ws = $wb.Worksheets.Item('East Village')
$ws.Range('J3').Value = 1
$ws.Range('J6').Value = 5
$ws.Range('J7').Value = 8

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J2').Value = 23
$ws.Range('J7').Value = 71

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J4').Value = 6
$ws.Range('J7').Value = 54

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('J3').Value = 17
$ws.Range('J7').Value = 60

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 115
$ws.Range('J3').Value = 121
$ws.Range('J6').Value = 96
$ws.Range('J7').Value = 359

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('J3').Value = 2
$ws.Range('J7').Value = 16

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('J4').Value = 14
$ws.Range('J7').Value = 31

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('J6').Value = 12
$ws.Range('J7').Value = 33

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J3').Value = 14
$ws.Range('J7').Value = 54

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('J3').Value = 6
$ws.Range('J7').Value = 28

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J3').Value = 13
$ws.Range('J7').Value = 43

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('J6').Value = 8
$ws.Range('J7').Value = 14

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J2').Value = 55
$ws.Range('J7').Value = 165
